# add exception handling for file upload
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Old town white coffee"
$ws.Range("B3").Value = "red bag"
$ws.Range("B4").Value = "metal plate"

$ws.Range("B5").Select()
